# Apply the "mixing" parameter rename/update to the "constant" sheet.
#
# Summary of the change (per the target diff):
#   - Row 3 parameter name: "child_socialising"   -> "mixing_a_spread"
#   - Row 4 parameter name: "elderly_socialising"  -> "mixing_pc_strength"
#   - Row 3: value (B3) 1 -> 10 ; distri_param1 (D3) 0.2 -> 5 ; distri_param2 (E3) 1 -> 20
#   - Row 4: distri_param1 (D4) 0.2 -> 1 ; distri_param2 (E4) 1 -> 5  (B4 unchanged)
#   - The "full_text" column (G) entries for rows 3 and 4 are removed entirely,
#     since the new parameters no longer carry a full_text description.
#   - Selection/active cell on the sheet moves to G4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constant")

# Rename the parameters in column A.
$ws.Range("A3").Value = "mixing_a_spread"
$ws.Range("A4").Value = "mixing_pc_strength"

# Update the numeric parameters for row 3 (mixing_a_spread).
$ws.Range("B3").Value = 10
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 20

# Update the numeric parameters for row 4 (mixing_pc_strength).
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 5

# Drop the now-unused "full_text" descriptions for these two rows.
$ws.Range("G3").Clear()
$ws.Range("G4").Clear()

# Move the active selection to reflect where the author's cursor ended up.
$ws.Range("G4").Select()
